$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 37/38 swap places (Hedera now ranks above TrustWalletToken) plus value updates
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05606"
$ws.Range("E37").Value = "  -3.59%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.142"
$ws.Range("E38").Value = "  -2.99%  "

# Price (D) and Volume(1h) (E) updates for the remaining rows
# D values that parse as plain numbers get a leading apostrophe so they are
# stored as text (matching the inlineStr text in the source file) instead of
# being auto-converted to numeric values by Excel.
$ws.Range("D2").Value = "28.659.44"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "1.886.53"
$ws.Range("E3").Value = "  -3.59%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'321.60"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.4562"
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("D8").Value = "'0.3786"
$ws.Range("E8").Value = "  -3.82%  "
$ws.Range("D9").Value = "'45.49"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").Value = "'0.07675"
$ws.Range("E10").Value = "  -2.73%  "
$ws.Range("D11").Value = "'0.9580"
$ws.Range("E11").Value = "  -4.31%  "
$ws.Range("D12").Value = "'21.89"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").Value = "1.891.94"
$ws.Range("E13").Value = "  -3.96%  "
$ws.Range("D14").Value = "'6.942"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "'5.632"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("D16").Value = "'0.07023"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "'82.49"
$ws.Range("E18").Value = "  -6.70%  "
$ws.Range("D19").Value = "'0.000009470"
$ws.Range("E19").Value = "  -4.90%  "
$ws.Range("D20").Value = "'16.58"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "28.645.04"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").Value = "'5.315"
$ws.Range("E23").Value = "  -3.74%  "
$ws.Range("D24").Value = "'10.81"
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").Value = "2.118.98"
$ws.Range("E25").Value = "  -3.98%  "
$ws.Range("D26").Value = "'2.059"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("D27").Value = "'155.05"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").Value = "'18.95"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("D29").Value = "'5.585"
$ws.Range("E29").Value = "  -6.60%  "
$ws.Range("D30").Value = "'116.49"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("D31").Value = "'1.808"
$ws.Range("E31").Value = "  -4.03%  "
$ws.Range("D32").Value = "'0.09202"
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("D33").Value = "'0.8420"
$ws.Range("E33").Value = "  -6.09%  "
$ws.Range("D34").Value = "'5.040"
$ws.Range("E34").Value = "  -3.88%  "
$ws.Range("D35").Value = "'1.241"
$ws.Range("E35").Value = "  -7.11%  "
$ws.Range("D36").Value = "'3.050"
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("D39").Value = "'1.001"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'0.02024"
$ws.Range("E40").Value = "  -4.52%  "
$ws.Range("D41").Value = "'7.419"
$ws.Range("E41").Value = "  -5.03%  "
$ws.Range("D42").Value = "'0.5459"
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("D43").Value = "'0.1740"
$ws.Range("E43").Value = "  -4.58%  "
$ws.Range("D44").Value = "'0.000002961"
$ws.Range("E44").Value = "  -22.34%  "
$ws.Range("D45").Value = "'9.128"
$ws.Range("E45").Value = "  -6.77%  "
$ws.Range("D46").Value = "'2.687"
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").Value = "'0.5141"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").Value = "'11.14"
$ws.Range("E48").Value = "  -8.26%  "
$ws.Range("D49").Value = "'2.083"
$ws.Range("E49").Value = "  -5.48%  "
$ws.Range("D50").Value = "'0.06753"
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").Value = "'110.40"
$ws.Range("E51").Value = "  -3.17%  "
